$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-04 Sunday" "2024-08-05 Monday"

Replace-Text "864÷6=144, 0" "810÷7=115, 5"
Replace-Text "252÷9=28, 0" "447÷8=55, 7"
Replace-Text "826÷9=91, 7" "713÷7=101, 6"
Replace-Text "390÷5=78, 0" "874÷5=174, 4"
Replace-Text "741÷7=105, 6" "101÷5=20, 1"

Replace-Text "532÷3=177, 1" "621÷7=88, 5"
Replace-Text "803÷8=100, 3" "622÷5=124, 2"
Replace-Text "945÷7=135, 0" "663÷7=94, 5"
Replace-Text "356÷2=178, 0" "482÷4=120, 2"
Replace-Text "930÷5=186, 0" "915÷2=457, 1"

Replace-Text "127÷9=14, 1" "482÷3=160, 2"
Replace-Text "746÷7=106, 4" "528÷6=88, 0"
Replace-Text "848÷2=424, 0" "331÷9=36, 7"
Replace-Text "736÷6=122, 4" "647÷8=80, 7"
Replace-Text "712÷4=178, 0" "981÷4=245, 1"

Replace-Text "655÷6=109, 1" "880÷5=176, 0"
Replace-Text "594÷3=198, 0" "722÷3=240, 2"
Replace-Text "351÷7=50, 1" "519÷8=64, 7"
Replace-Text "519÷7=74, 1" "477÷7=68, 1"
Replace-Text "782÷3=260, 2" "632÷5=126, 2"

Replace-Text "114÷2=57, 0" "977÷8=122, 1"
Replace-Text "537÷7=76, 5" "984÷4=246, 0"
Replace-Text "819÷4=204, 3" "459÷6=76, 3"
Replace-Text "231÷3=77, 0" "480÷9=53, 3"
Replace-Text "873÷6=145, 3" "237÷3=79, 0"
